# Dynamic GW/SW change for Lake Annie.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# "DOC_gw = InLakeDOC" flag (row 27): turn the dynamic DOC_gw switch off.
$ws.Range("B27").Value = 0

# "ProdEndDay" (row 28): extend the production end day to cover a leap year.
$ws.Range("B28").Value = 366

# Reflect the user's final selection/scroll position on the sheet.
[void]$ws.Range("D28").Select()
